# Update column F (dSF) values for the rows that were repulled/recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -7
    4  = -1
    25 = -1
    26 = -3
    31 = -3
    36 = 5
    43 = -2
    45 = 2
    48 = -3
    50 = 2
    52 = 0
    53 = 4
    55 = -4
    56 = 2
    63 = -9
    64 = -2
    68 = -1
    72 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
